$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.015475869178772
$ws.Range("B1").Value = 1.506506323814392
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 0.438677579164505
